# Update the "dSF" (F) column values for a handful of rows in the sheet,
# reflecting a repull of the underlying data / recalculated mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    5  = 3
    6  = 4
    9  = -5
    10 = -1
    12 = -1
    18 = -1
    20 = 1
    24 = -1
    25 = -2
    26 = 2
    30 = -1
    31 = -1
    37 = -5
    40 = 6
    41 = 2
    43 = 8
    44 = 3
    46 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
